$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text edits (Volume number + report date range) ---
$ws.Range("A8").Value = "Volume 30   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/4/2023  Through  12/10/2023"

# --- Cells changing between numeric <-> text ("0"/"***.*") representation ---
# Use .Copy() from a stable style/type donor in row 14 (untouched by this edit)
# to get the right style + type, then overwrite numeric cells with their value.
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("I14").Copy($ws.Range("C17"))
$ws.Range("C14").Copy($ws.Range("C20"))
$ws.Range("I14").Copy($ws.Range("C22"))
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("C14").Copy($ws.Range("C26"))
$ws.Range("I14").Copy($ws.Range("C27"))
$ws.Range("I14").Copy($ws.Range("C28"))
$ws.Range("I14").Copy($ws.Range("F28"))
$ws.Range("I14").Copy($ws.Range("C29"))
$ws.Range("I14").Copy($ws.Range("F29"))
$ws.Range("C14").Copy($ws.Range("D30"))
$ws.Range("E14").Copy($ws.Range("E30"))

# Now set final numeric values for the cells we just re-typed via Copy
$ws.Range("C17").Value = 1
$ws.Range("C22").Value = 1
$ws.Range("C27").Value = 1
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 1

# --- Plain value-only updates (style/type unchanged) ---
$ws.Range("C16").Value = 8
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 69.230769230769
$ws.Range("I16").Value = 205
$ws.Range("J16").Value = 226
$ws.Range("K16").Value = -9.292035398230
$ws.Range("L16").Value = 21.301775147929
$ws.Range("M16").Value = 51.851851851851
$ws.Range("N16").Value = -84.279141104294
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -83.333333333333
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -61.904761904761
$ws.Range("I17").Value = 173
$ws.Range("J17").Value = 171
$ws.Range("K17").Value = 1.169590643274
$ws.Range("L17").Value = 16.107382550335
$ws.Range("M17").Value = 92.222222222222
$ws.Range("N17").Value = -38.434163701067
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 20
$ws.Range("E18").Value = -90
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 42
$ws.Range("H18").Value = -71.428571428571
$ws.Range("I18").Value = 243
$ws.Range("J18").Value = 271
$ws.Range("K18").Value = -10.332103321033
$ws.Range("L18").Value = 11.467889908256
$ws.Range("M18").Value = 7.048458149779
$ws.Range("N18").Value = -91.470691470691
$ws.Range("C19").Value = 42
$ws.Range("D19").Value = 37
$ws.Range("E19").Value = 13.513513513513
$ws.Range("F19").Value = 142
$ws.Range("G19").Value = 122
$ws.Range("H19").Value = 16.393442622950
$ws.Range("I19").Value = 1618
$ws.Range("J19").Value = 1659
$ws.Range("K19").Value = -2.471368294153
$ws.Range("L19").Value = 41.805433829973
$ws.Range("M19").Value = 35.397489539749
$ws.Range("N19").Value = -54.486638537271
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = -69.230769230769
$ws.Range("J20").Value = 181
$ws.Range("K20").Value = -18.784530386740
$ws.Range("M20").Value = 68.965517241379
$ws.Range("N20").Value = -95.525114155251
$ws.Range("C21").Value = 53
$ws.Range("D21").Value = 69
$ws.Range("E21").Value = -23.188405797101
$ws.Range("F21").Value = 192
$ws.Range("G21").Value = 225
$ws.Range("H21").Value = -14.666666666666
$ws.Range("I21").Value = 2401
$ws.Range("J21").Value = 2521
$ws.Range("K21").Value = -4.760015866719
$ws.Range("L21").Value = 29.853975121687
$ws.Range("M21").Value = 37.278444825614
$ws.Range("N21").Value = -78.769121938279
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 42
$ws.Range("K22").Value = -6.666666666666
$ws.Range("L22").Value = 133.333333333333
$ws.Range("M22").Value = 68
$ws.Range("F23").Value = 2
$ws.Range("H23").Value = -66.666666666666
$ws.Range("J23").Value = 38
$ws.Range("K23").Value = -26.315789473684
$ws.Range("L23").Value = -9.677419354838
$ws.Range("C24").Value = 47
$ws.Range("D24").Value = 86
$ws.Range("E24").Value = -45.348837209302
$ws.Range("F24").Value = 206
$ws.Range("G24").Value = 348
$ws.Range("H24").Value = -40.804597701149
$ws.Range("I24").Value = 2935
$ws.Range("J24").Value = 3835
$ws.Range("K24").Value = -23.468057366362
$ws.Range("L24").Value = 25.427350427350
$ws.Range("M24").Value = 84.70736312146
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -28.571428571428
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 28
$ws.Range("H25").Value = 35.714285714285
$ws.Range("I25").Value = 339
$ws.Range("J25").Value = 356
$ws.Range("K25").Value = -4.775280898876
$ws.Range("L25").Value = 5.607476635514
$ws.Range("M25").Value = 3.039513677811
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -75
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = 33.333333333333
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 82
$ws.Range("J27").Value = 81
$ws.Range("K27").Value = 1.234567901234
$ws.Range("L27").Value = -6.818181818181
$ws.Range("I28").Value = 4
$ws.Range("K28").Value = 100
$ws.Range("L28").Value = -42.857142857142
$ws.Range("M28").Value = 33.333333333333
$ws.Range("N28").Value = 100
$ws.Range("I29").Value = 4
$ws.Range("K29").Value = 100
$ws.Range("L29").Value = -33.333333333333
$ws.Range("M29").Value = 33.333333333333
$ws.Range("N29").Value = 100
$ws.Range("F30").Value = 4
$ws.Range("H30").Value = 300
$ws.Range("I30").Value = 23
$ws.Range("K30").Value = -17.857142857142
$ws.Range("L30").Value = 53.333333333333
